# Weekly update: a new day's price record is added for Albahaca at
# "Terminal La Palmera de La Serena" (Coquimbo). This inserts a new row
# at row 9 (pushing the existing rows 9-81 down to 10-82) and populates
# the new row with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing data rows down by one.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data.
$ws.Cells.Item(9, 1).Value  = 8
$ws.Cells.Item(9, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(9, 3).Value  = "Coquimbo"
$ws.Cells.Item(9, 4).Value  = 44670
$ws.Cells.Item(9, 5).Value  = 4
$ws.Cells.Item(9, 6).Value  = 100112052
$ws.Cells.Item(9, 7).Value  = "Albahaca"
$ws.Cells.Item(9, 8).Value  = "Sin especificar"
$ws.Cells.Item(9, 9).Value  = "Primera"
$ws.Cells.Item(9, 10).Value = 800
$ws.Cells.Item(9, 11).Value = 5000
$ws.Cells.Item(9, 12).Value = 5500
$ws.Cells.Item(9, 13).Value = 5250
$ws.Cells.Item(9, 14).Value = "$/docena de matas"
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(9, 16).Value = 875
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = "Hortaliza"
